$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.80"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.14"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.347"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05935"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.396"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.394"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9662"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1431"

$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07395"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"

$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03489"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03045"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09407"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.002"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001601"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04804"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005913"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006173"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.004144"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0009839"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.00009706"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.743"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03936"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006484"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002702"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.005376"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005304"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.8504"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.04433"
